$d = $word.ActiveDocument
Write-Host "=== FULL TEXT ==="
Write-Host $d.Content.Text
